# Update column C ("Förändrad") for every data row: the date serial
# 45203 (2023-10-04) becomes 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCell = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162)
$lastRow = $lastCell.Row()
if ($lastRow -lt 2) { $lastRow = 382 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 45203) {
        $cell.Value = 45204
    }
}
